# Apply the 2024-05-21 update to slide 1 of the presentation:
#  1. Update the presentation date from 2024.05.20 to 2024.05.21
#  2. Update the presenter name from 오규안 to 김태환 (splitting the
#     "발표자 : 오규안" line into separate runs for the label, the
#     colon separator and the name, matching the authored run layout)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- 1. Date line -------------------------------------------------------
$dateParagraph = $tr.Paragraphs(2)
$dateChars = $tr.Characters($dateParagraph.Start, $dateParagraph.Length)
$dateChars.Text = "2024.05.21"

# --- 2. Presenter line ---------------------------------------------------
$presenterParagraph = $tr.Paragraphs(5)
$pStart = $presenterParagraph.Start

# "발표자 : 오규안"
#   chars 1-4 => "발표자 "   (label, keeps its own run/formatting)
#   chars 5-6 => ": "        (separator becomes its own run)
#   chars 7-9 => "오규안"    (name -> replaced with "김태환")
$label = $tr.Characters($pStart, 4)
$label.Text = "발표자 "

$colon = $tr.Characters($pStart + 4, 2)
$colon.Text = ": "

$name = $tr.Characters($pStart + 6, 3)
$name.Text = "김태환"
